$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{R=14; C=2; V=6838331},
    @{R=14; C=5; V='Feyenoord'},
    @{R=14; C=6; V='Fortuna Sittard'},
    @{R=14; C=7; V=0},
    @{R=14; C=8; V=0},
    @{R=14; C=9; V=0},
    @{R=14; C=11; V='D'},
    @{R=14; C=12; V=1.25},
    @{R=14; C=13; V=6},
    @{R=14; C=14; V=10},
    @{R=14; C=15; V=1.2},
    @{R=14; C=16; V=7},
    @{R=14; C=17; V=17},
    @{R=14; C=18; V=-2.25},
    @{R=14; C=19; V=2.05},
    @{R=14; C=20; V=1.8},
    @{R=14; C=21; V=3.5},
    @{R=14; C=22; V=2},
    @{R=14; C=23; V=1.85},
    @{R=14; C=24; V=-1},
    @{R=14; C=25; V=6},
    @{R=14; C=27; V=-1},
    @{R=14; C=28; V=0.8},
    @{R=14; C=29; V=-1},
    @{R=14; C=30; V=0.8500000000000001},
    @{R=15; C=2; V=6838332},
    @{R=15; C=5; V='AZ'},
    @{R=15; C=6; V='Go Ahead Eagles'},
    @{R=15; C=7; V=5},
    @{R=15; C=8; V=1},
    @{R=15; C=9; V=2},
    @{R=15; C=11; V='H'},
    @{R=15; C=12; V=1.4},
    @{R=15; C=13; V=4.5},
    @{R=15; C=14; V=8},
    @{R=15; C=15; V=1.55},
    @{R=15; C=16; V=4},
    @{R=15; C=17; V=6},
    @{R=15; C=18; V=-1},
    @{R=15; C=19; V=1.9},
    @{R=15; C=20; V=1.95},
    @{R=15; C=21; V=2.75},
    @{R=15; C=22; V=1.95},
    @{R=15; C=23; V=1.9},
    @{R=15; C=24; V=0.55},
    @{R=15; C=25; V=-1},
    @{R=15; C=27; V=0.8999999999999999},
    @{R=15; C=28; V=-1},
    @{R=15; C=29; V=0.95},
    @{R=15; C=30; V=-1},
    @{R=88; C=2; V=6838413},
    @{R=88; C=5; V='Sparta Rotterdam'},
    @{R=88; C=6; V='RKC'},
    @{R=88; C=7; V=2},
    @{R=88; C=11; V='H'},
    @{R=88; C=12; V=1.8},
    @{R=88; C=13; V=3.6},
    @{R=88; C=14; V=4.333},
    @{R=88; C=15; V=1.75},
    @{R=88; C=16; V=3.75},
    @{R=88; C=17; V=4.75},
    @{R=88; C=18; V=-0.75},
    @{R=88; C=19; V=1.975},
    @{R=88; C=20; V=1.875},
    @{R=88; C=22; V=1.925},
    @{R=88; C=23; V=1.925},
    @{R=88; C=24; V=0.75},
    @{R=88; C=25; V=-1},
    @{R=88; C=27; V=0.9750000000000001},
    @{R=88; C=28; V=-1},
    @{R=88; C=30; V=0.925},
    @{R=89; C=2; V=6838411},
    @{R=89; C=5; V='Almere City FC'},
    @{R=89; C=6; V='Go Ahead Eagles'},
    @{R=89; C=7; V=0},
    @{R=89; C=11; V='D'},
    @{R=89; C=12; V=2.625},
    @{R=89; C=13; V=3.5},
    @{R=89; C=14; V=2.5},
    @{R=89; C=15; V=2.9},
    @{R=89; C=16; V=3.6},
    @{R=89; C=17; V=2.3},
    @{R=89; C=18; V=0.25},
    @{R=89; C=19; V=1.875},
    @{R=89; C=20; V=1.975},
    @{R=89; C=22; V=1.9},
    @{R=89; C=23; V=1.95},
    @{R=89; C=24; V=-1},
    @{R=89; C=25; V=2.6},
    @{R=89; C=27; V=0.4375},
    @{R=89; C=28; V=-0.5},
    @{R=89; C=30; V=0.95},
    @{R=102; C=2; V=6838426},
    @{R=102; C=5; V='Sparta Rotterdam'},
    @{R=102; C=6; V='Almere City FC'},
    @{R=102; C=7; V=1},
    @{R=102; C=8; V=2},
    @{R=102; C=9; V=0},
    @{R=102; C=11; V='A'},
    @{R=102; C=12; V=2.15},
    @{R=102; C=13; V=3.4},
    @{R=102; C=14; V=3},
    @{R=102; C=15; V=1.7},
    @{R=102; C=16; V=4},
    @{R=102; C=17; V=4.75},
    @{R=102; C=18; V=-0.75},
    @{R=102; C=19; V=1.9},
    @{R=102; C=20; V=1.95},
    @{R=102; C=21; V=2.75},
    @{R=102; C=22; V=2.025},
    @{R=102; C=23; V=1.825},
    @{R=102; C=24; V=-1},
    @{R=102; C=26; V=3.75},
    @{R=102; C=27; V=-1},
    @{R=102; C=28; V=0.95},
    @{R=102; C=29; V=0.5125},
    @{R=102; C=30; V=-0.5},
    @{R=103; C=2; V=6838419},
    @{R=103; C=5; V='Ajax'},
    @{R=103; C=6; V='Heerenveen'},
    @{R=103; C=7; V=4},
    @{R=103; C=8; V=1},
    @{R=103; C=9; V=2},
    @{R=103; C=11; V='H'},
    @{R=103; C=12; V=1.444},
    @{R=103; C=13; V=4.6},
    @{R=103; C=14; V=5.5},
    @{R=103; C=15; V=1.4},
    @{R=103; C=16; V=5.25},
    @{R=103; C=17; V=7},
    @{R=103; C=18; V=-1.5},
    @{R=103; C=19; V=2},
    @{R=103; C=20; V=1.85},
    @{R=103; C=21; V=3.5},
    @{R=103; C=22; V=1.975},
    @{R=103; C=23; V=1.875},
    @{R=103; C=24; V=0.3999999999999999},
    @{R=103; C=26; V=-1},
    @{R=103; C=27; V=1},
    @{R=103; C=28; V=-1},
    @{R=103; C=29; V=0.9750000000000001},
    @{R=103; C=30; V=-1},
    @{R=114; C=2; V=6838437},
    @{R=114; C=5; V='FC Twente'},
    @{R=114; C=6; V='PSV'},
    @{R=114; C=7; V=0},
    @{R=114; C=8; V=3},
    @{R=114; C=9; V=0},
    @{R=114; C=10; V=1},
    @{R=114; C=11; V='A'},
    @{R=114; C=12; V=4},
    @{R=114; C=14; V=1.833},
    @{R=114; C=15; V=3.6},
    @{R=114; C=16; V=3.6},
    @{R=114; C=17; V=2},
    @{R=114; C=18; V=0.5},
    @{R=114; C=19; V=1.825},
    @{R=114; C=20; V=2.025},
    @{R=114; C=21; V=3},
    @{R=114; C=22; V=2.025},
    @{R=114; C=23; V=1.825},
    @{R=114; C=24; V=-1},
    @{R=114; C=26; V=1},
    @{R=114; C=27; V=-1},
    @{R=114; C=28; V=1.025},
    @{R=114; C=29; V=0},
    @{R=114; C=30; V=0},
    @{R=115; C=2; V=6838438},
    @{R=115; C=5; V='Heerenveen'},
    @{R=115; C=6; V='Fortuna Sittard'},
    @{R=115; C=7; V=3},
    @{R=115; C=8; V=0},
    @{R=115; C=9; V=1},
    @{R=115; C=10; V=0},
    @{R=115; C=11; V='H'},
    @{R=115; C=12; V=1.833},
    @{R=115; C=14; V=4},
    @{R=115; C=15; V=1.95},
    @{R=115; C=16; V=3.5},
    @{R=115; C=17; V=3.8},
    @{R=115; C=18; V=-0.5},
    @{R=115; C=19; V=1.975},
    @{R=115; C=20; V=1.875},
    @{R=115; C=21; V=2.5},
    @{R=115; C=22; V=2},
    @{R=115; C=23; V=1.85},
    @{R=115; C=24; V=0.95},
    @{R=115; C=26; V=-1},
    @{R=115; C=27; V=0.9750000000000001},
    @{R=115; C=28; V=-1},
    @{R=115; C=29; V=1},
    @{R=115; C=30; V=-1},
    @{R=116; C=2; V=6838440},
    @{R=116; C=5; V='Ajax'},
    @{R=116; C=6; V='Vitesse'},
    @{R=116; C=7; V=5},
    @{R=116; C=8; V=0},
    @{R=116; C=9; V=2},
    @{R=116; C=10; V=0},
    @{R=116; C=11; V='H'},
    @{R=116; C=12; V=1.363},
    @{R=116; C=13; V=5},
    @{R=116; C=14; V=7.5},
    @{R=116; C=15; V=1.333},
    @{R=116; C=16; V=6},
    @{R=116; C=17; V=7.5},
    @{R=116; C=18; V=-1.5},
    @{R=116; C=19; V=1.825},
    @{R=116; C=20; V=2.025},
    @{R=116; C=21; V=3.5},
    @{R=116; C=22; V=1.95},
    @{R=116; C=23; V=1.9},
    @{R=116; C=24; V=0.333},
    @{R=116; C=26; V=-1},
    @{R=116; C=27; V=0.825},
    @{R=116; C=28; V=-1},
    @{R=116; C=29; V=0.95},
    @{R=116; C=30; V=-1},
    @{R=117; C=2; V=6838439},
    @{R=117; C=5; V='PEC Zwolle'},
    @{R=117; C=6; V='RKC'},
    @{R=117; C=7; V=1},
    @{R=117; C=8; V=2},
    @{R=117; C=9; V=1},
    @{R=117; C=10; V=1},
    @{R=117; C=11; V='A'},
    @{R=117; C=12; V=1.909},
    @{R=117; C=13; V=3.6},
    @{R=117; C=14; V=3.75},
    @{R=117; C=15; V=2.05},
    @{R=117; C=16; V=3.6},
    @{R=117; C=17; V=3.4},
    @{R=117; C=18; V=-0.5},
    @{R=117; C=19; V=2.05},
    @{R=117; C=20; V=1.8},
    @{R=117; C=21; V=2.75},
    @{R=117; C=22; V=1.9},
    @{R=117; C=23; V=1.95},
    @{R=117; C=24; V=-1},
    @{R=117; C=26; V=2.4},
    @{R=117; C=27; V=-1},
    @{R=117; C=28; V=0.8},
    @{R=117; C=29; V=0.45},
    @{R=117; C=30; V=-0.5},
    @{R=124; C=2; V=6838448},
    @{R=124; C=5; V='RKC'},
    @{R=124; C=6; V='Excelsior'},
    @{R=124; C=7; V=2},
    @{R=124; C=8; V=2},
    @{R=124; C=9; V=1},
    @{R=124; C=11; V='D'},
    @{R=124; C=12; V=2.05},
    @{R=124; C=13; V=3.7},
    @{R=124; C=14; V=3.1},
    @{R=124; C=15; V=1.95},
    @{R=124; C=16; V=3.8},
    @{R=124; C=17; V=3.5},
    @{R=124; C=19; V=2},
    @{R=124; C=20; V=1.85},
    @{R=124; C=21; V=2.75},
    @{R=124; C=22; V=1.85},
    @{R=124; C=23; V=2},
    @{R=124; C=24; V=-1},
    @{R=124; C=25; V=2.8},
    @{R=124; C=27; V=-1},
    @{R=124; C=28; V=0.8500000000000001},
    @{R=124; C=29; V=0.8500000000000001},
    @{R=125; C=2; V=6838447},
    @{R=125; C=5; V='Fortuna Sittard'},
    @{R=125; C=6; V='Vitesse'},
    @{R=125; C=7; V=3},
    @{R=125; C=8; V=1},
    @{R=125; C=9; V=2},
    @{R=125; C=11; V='H'},
    @{R=125; C=12; V=2},
    @{R=125; C=13; V=3.5},
    @{R=125; C=14; V=3.4},
    @{R=125; C=15; V=2},
    @{R=125; C=16; V=3.6},
    @{R=125; C=17; V=3.6},
    @{R=125; C=19; V=2.025},
    @{R=125; C=20; V=1.825},
    @{R=125; C=21; V=2.5},
    @{R=125; C=22; V=1.875},
    @{R=125; C=23; V=1.975},
    @{R=125; C=24; V=1},
    @{R=125; C=25; V=-1},
    @{R=125; C=27; V=1.025},
    @{R=125; C=28; V=-1},
    @{R=125; C=29; V=0.875},
    @{R=155; C=2; V=6838476},
    @{R=155; C=5; V='Go Ahead Eagles'},
    @{R=155; C=6; V='Ajax'},
    @{R=155; C=7; V=2},
    @{R=155; C=8; V=3},
    @{R=155; C=9; V=1},
    @{R=155; C=10; V=2},
    @{R=155; C=11; V='A'},
    @{R=155; C=12; V=3.5},
    @{R=155; C=14; V=2},
    @{R=155; C=15; V=4.2},
    @{R=155; C=16; V=4},
    @{R=155; C=17; V=1.75},
    @{R=155; C=18; V=0.75},
    @{R=155; C=19; V=1.85},
    @{R=155; C=20; V=2},
    @{R=155; C=21; V=3.25},
    @{R=155; C=22; V=1.975},
    @{R=155; C=23; V=1.875},
    @{R=155; C=25; V=-1},
    @{R=155; C=26; V=0.75},
    @{R=155; C=27; V=-0.5},
    @{R=155; C=28; V=0.5},
    @{R=155; C=29; V=0.9750000000000001},
    @{R=155; C=30; V=-1},
    @{R=156; C=2; V=6838477},
    @{R=156; C=5; V='Vitesse'},
    @{R=156; C=6; V='FC Utrecht'},
    @{R=156; C=7; V=0},
    @{R=156; C=8; V=0},
    @{R=156; C=9; V=0},
    @{R=156; C=10; V=0},
    @{R=156; C=11; V='D'},
    @{R=156; C=12; V=2.6},
    @{R=156; C=14; V=2.5},
    @{R=156; C=15; V=2.9},
    @{R=156; C=16; V=3.8},
    @{R=156; C=17; V=2.2},
    @{R=156; C=18; V=0.25},
    @{R=156; C=19; V=1.825},
    @{R=156; C=20; V=2.025},
    @{R=156; C=21; V=2.75},
    @{R=156; C=22; V=2},
    @{R=156; C=23; V=1.85},
    @{R=156; C=25; V=2.8},
    @{R=156; C=26; V=-1},
    @{R=156; C=27; V=0.4125},
    @{R=156; C=28; V=-0.5},
    @{R=156; C=29; V=-1},
    @{R=156; C=30; V=0.8500000000000001},
    @{R=190; C=2; V=6838506},
    @{R=190; C=5; V='FC Utrecht'},
    @{R=190; C=6; V='Fortuna Sittard'},
    @{R=190; C=7; V=4},
    @{R=190; C=8; V=0},
    @{R=190; C=9; V=1},
    @{R=190; C=12; V=2.1},
    @{R=190; C=13; V=3.3},
    @{R=190; C=14; V=3.5},
    @{R=190; C=15; V=1.615},
    @{R=190; C=16; V=4},
    @{R=190; C=17; V=5.5},
    @{R=190; C=18; V=-0.75},
    @{R=190; C=19; V=1.8},
    @{R=190; C=20; V=2.05},
    @{R=190; C=21; V=2.5},
    @{R=190; C=22; V=1.95},
    @{R=190; C=23; V=1.9},
    @{R=190; C=24; V=0.615},
    @{R=190; C=27; V=0.8},
    @{R=190; C=29; V=0.95},
    @{R=191; C=2; V=6838511},
    @{R=191; C=5; V='Heerenveen'},
    @{R=191; C=6; V='Ajax'},
    @{R=191; C=7; V=3},
    @{R=191; C=8; V=2},
    @{R=191; C=9; V=2},
    @{R=191; C=12; V=5},
    @{R=191; C=13; V=4.333},
    @{R=191; C=14; V=1.571},
    @{R=191; C=15; V=3.5},
    @{R=191; C=16; V=4.333},
    @{R=191; C=17; V=1.85},
    @{R=191; C=18; V=0.5},
    @{R=191; C=19; V=1.925},
    @{R=191; C=20; V=1.925},
    @{R=191; C=21; V=3.25},
    @{R=191; C=22; V=1.9},
    @{R=191; C=23; V=1.95},
    @{R=191; C=24; V=2.5},
    @{R=191; C=27; V=0.925},
    @{R=191; C=29; V=0.8999999999999999},
    @{R=208; C=2; V=6838524},
    @{R=208; C=5; V='FC Twente'},
    @{R=208; C=6; V='Go Ahead Eagles'},
    @{R=208; C=7; V=3},
    @{R=208; C=8; V=0},
    @{R=208; C=9; V=2},
    @{R=208; C=11; V='H'},
    @{R=208; C=12; V=1.4},
    @{R=208; C=13; V=4.8},
    @{R=208; C=14; V=7},
    @{R=208; C=15; V=1.363},
    @{R=208; C=16; V=5.25},
    @{R=208; C=17; V=7.5},
    @{R=208; C=18; V=-1.5},
    @{R=208; C=19; V=2.025},
    @{R=208; C=20; V=1.825},
    @{R=208; C=21; V=3},
    @{R=208; C=22; V=1.95},
    @{R=208; C=23; V=1.9},
    @{R=208; C=24; V=0.363},
    @{R=208; C=26; V=-1},
    @{R=208; C=27; V=1.025},
    @{R=208; C=28; V=-1},
    @{R=208; C=29; V=0},
    @{R=208; C=30; V=0},
    @{R=209; C=2; V=6838521},
    @{R=209; C=5; V='Almere City FC'},
    @{R=209; C=6; V='Feyenoord'},
    @{R=209; C=7; V=0},
    @{R=209; C=8; V=2},
    @{R=209; C=9; V=0},
    @{R=209; C=11; V='A'},
    @{R=209; C=12; V=8.5},
    @{R=209; C=13; V=5.6},
    @{R=209; C=14; V=1.3},
    @{R=209; C=15; V=12},
    @{R=209; C=16; V=5.5},
    @{R=209; C=17; V=1.25},
    @{R=209; C=18; V=1.75},
    @{R=209; C=19; V=1.825},
    @{R=209; C=20; V=2.025},
    @{R=209; C=21; V=2.75},
    @{R=209; C=22; V=1.825},
    @{R=209; C=23; V=2.025},
    @{R=209; C=24; V=-1},
    @{R=209; C=26; V=0.25},
    @{R=209; C=27; V=-0.5},
    @{R=209; C=28; V=0.5125},
    @{R=209; C=29; V=-1},
    @{R=209; C=30; V=1.025},
    @{R=235; C=2; V=6838548},
    @{R=235; C=5; V='FC Utrecht'},
    @{R=235; C=6; V='NEC'},
    @{R=235; C=7; V=1},
    @{R=235; C=8; V=0},
    @{R=235; C=9; V=0},
    @{R=235; C=10; V=0},
    @{R=235; C=11; V='H'},
    @{R=235; C=12; V=2.4},
    @{R=235; C=13; V=3.25},
    @{R=235; C=14; V=3},
    @{R=235; C=15; V=2.05},
    @{R=235; C=16; V=3.4},
    @{R=235; C=17; V=3.6},
    @{R=235; C=18; V=-0.5},
    @{R=235; C=19; V=2.07},
    @{R=235; C=20; V=1.83},
    @{R=235; C=21; V=2.5},
    @{R=235; C=22; V=1.825},
    @{R=235; C=23; V=2.025},
    @{R=235; C=24; V=1.05},
    @{R=235; C=26; V=-1},
    @{R=235; C=27; V=1.07},
    @{R=235; C=28; V=-1},
    @{R=235; C=29; V=-1},
    @{R=235; C=30; V=1.025},
    @{R=236; C=2; V=6838552},
    @{R=236; C=5; V='Heerenveen'},
    @{R=236; C=6; V='Feyenoord'},
    @{R=236; C=7; V=2},
    @{R=236; C=8; V=3},
    @{R=236; C=9; V=1},
    @{R=236; C=10; V=1},
    @{R=236; C=11; V='A'},
    @{R=236; C=12; V=4.333},
    @{R=236; C=13; V=4.5},
    @{R=236; C=14; V=1.65},
    @{R=236; C=15; V=6.5},
    @{R=236; C=16; V=5.25},
    @{R=236; C=17; V=1.4},
    @{R=236; C=18; V=1.25},
    @{R=236; C=19; V=2},
    @{R=236; C=20; V=1.85},
    @{R=236; C=21; V=3},
    @{R=236; C=22; V=2.025},
    @{R=236; C=23; V=1.825},
    @{R=236; C=24; V=-1},
    @{R=236; C=26; V=0.3999999999999999},
    @{R=236; C=27; V=0.5},
    @{R=236; C=28; V=-0.5},
    @{R=236; C=29; V=1.025},
    @{R=236; C=30; V=-1},
    @{R=302; C=2; V=7223358},
    @{R=302; C=5; V='Vitesse'},
    @{R=302; C=6; V='Ajax'},
    @{R=302; C=7; V=2},
    @{R=302; C=8; V=2},
    @{R=302; C=9; V=1},
    @{R=302; C=10; V=1},
    @{R=302; C=12; V=4.75},
    @{R=302; C=13; V=4},
    @{R=302; C=14; V=1.571},
    @{R=302; C=15; V=4.333},
    @{R=302; C=16; V=4.5},
    @{R=302; C=17; V=1.65},
    @{R=302; C=18; V=1},
    @{R=302; C=19; V=1.84},
    @{R=302; C=20; V=2.06},
    @{R=302; C=21; V=3.75},
    @{R=302; C=25; V=3.5},
    @{R=302; C=27; V=0.8400000000000001},
    @{R=302; C=28; V=-1},
    @{R=302; C=29; V=0.4625},
    @{R=302; C=30; V=-0.5},
    @{R=303; C=2; V=7160673},
    @{R=303; C=5; V='AZ'},
    @{R=303; C=6; V='FC Utrecht'},
    @{R=303; C=7; V=3},
    @{R=303; C=8; V=3},
    @{R=303; C=9; V=3},
    @{R=303; C=10; V=0},
    @{R=303; C=11; V='D'},
    @{R=303; C=12; V=1.4},
    @{R=303; C=13; V=4.5},
    @{R=303; C=14; V=6.5},
    @{R=303; C=15; V=1.4},
    @{R=303; C=16; V=5},
    @{R=303; C=17; V=7},
    @{R=303; C=18; V=-1.25},
    @{R=303; C=21; V=3.25},
    @{R=303; C=22; V=1.925},
    @{R=303; C=23; V=1.925},
    @{R=303; C=25; V=4},
    @{R=303; C=26; V=-1},
    @{R=303; C=29; V=0.925},
    @{R=304; C=2; V=7155056},
    @{R=304; C=5; V='Almere City FC'},
    @{R=304; C=6; V='NEC'},
    @{R=304; C=7; V=1},
    @{R=304; C=8; V=4},
    @{R=304; C=9; V=1},
    @{R=304; C=10; V=2},
    @{R=304; C=11; V='A'},
    @{R=304; C=12; V=3.1},
    @{R=304; C=13; V=3.4},
    @{R=304; C=14; V=2.1},
    @{R=304; C=15; V=3},
    @{R=304; C=16; V=3.6},
    @{R=304; C=17; V=2.3},
    @{R=304; C=18; V=0.25},
    @{R=304; C=19; V=1.825},
    @{R=304; C=20; V=2.025},
    @{R=304; C=21; V=2.75},
    @{R=304; C=22; V=1.8},
    @{R=304; C=23; V=2.05},
    @{R=304; C=25; V=-1},
    @{R=304; C=26; V=1.3},
    @{R=304; C=27; V=-1},
    @{R=304; C=28; V=1.025},
    @{R=304; C=29; V=0.8},
    @{R=304; C=30; V=-1},
    @{R=305; C=2; V=7161289},
    @{R=305; C=5; V='Heracles'},
    @{R=305; C=6; V='Fortuna Sittard'},
    @{R=305; C=7; V=0},
    @{R=305; C=8; V=0},
    @{R=305; C=9; V=0},
    @{R=305; C=10; V=0},
    @{R=305; C=11; V='D'},
    @{R=305; C=12; V=2.25},
    @{R=305; C=13; V=3.5},
    @{R=305; C=14; V=2.75},
    @{R=305; C=15; V=2.6},
    @{R=305; C=16; V=3.7},
    @{R=305; C=17; V=2.5},
    @{R=305; C=18; V=0},
    @{R=305; C=19; V=1.975},
    @{R=305; C=20; V=1.875},
    @{R=305; C=21; V=3.25},
    @{R=305; C=22; V=2.025},
    @{R=305; C=23; V=1.825},
    @{R=305; C=25; V=2.7},
    @{R=305; C=26; V=-1},
    @{R=305; C=30; V=0.825},
    @{R=306; C=2; V=7223357},
    @{R=306; C=5; V='FC Volendam'},
    @{R=306; C=6; V='Go Ahead Eagles'},
    @{R=306; C=7; V=1},
    @{R=306; C=11; V='A'},
    @{R=306; C=12; V=5},
    @{R=306; C=13; V=4.333},
    @{R=306; C=14; V=1.5},
    @{R=306; C=15; V=5.25},
    @{R=306; C=16; V=4.75},
    @{R=306; C=17; V=1.533},
    @{R=306; C=19; V=2.05},
    @{R=306; C=20; V=1.85},
    @{R=306; C=22; V=1.875},
    @{R=306; C=23; V=1.975},
    @{R=306; C=25; V=-1},
    @{R=306; C=26; V=0.5329999999999999},
    @{R=306; C=27; V=0},
    @{R=306; C=28; V=0},
    @{R=306; C=29; V=-1},
    @{R=306; C=30; V=0.9750000000000001}
)

foreach ($item in $changes) {
    $ws.Cells.Item($item.R, $item.C).Value = $item.V
}
